$wb = $excel.ActiveWorkbook

# The change affects both the "展览" sheet and the "全部类型" sheet, which
# contain duplicate rows of the same convention data ("想去人数" values).
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 2: F2 想去人数 9559 -> 9589
    $ws.Range("F2").Value = 9589

    # Row 4: F4 想去人数 28 -> 30
    $ws.Range("F4").Value = 30

    # Row 5: F5 想去人数 534 -> 540
    $ws.Range("F5").Value = 540
}
